$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 13-14: height 21 -> 20.25
# Rows 15-24: height 18.75 -> 20.25
$ws.Range("A13:A24").EntireRow.RowHeight = 20.25
